$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "REQ_ PO2_EBL_Electric_Blender_SRS_001.3-1.0"
$ws.Range("B2").Value = "Mostafa Ramadan"
$ws.Range("C2").Value = "what should be the speed for each state?"
